# Update column G ("K") values for rows 2-19 on the active sheet.
# These reflect a regen of the save_data using "K" (strikeouts) in place
# of the previous "Strike#" computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 4
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 4
    11 = 2
    12 = 1
    13 = 1
    14 = 3
    15 = 3
    16 = 5
    17 = 2
    18 = 3
    19 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
